$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 98, shifting existing rows
# (old 98..129) down to (100..131).
$ws.Rows.Item(98).Insert()
$ws.Rows.Item(98).Insert()

# Populate the two newly inserted rows (98 and 99) with the new weekly
# price entries for "Ají" (Inferno).

# Row 98: Inferno / Primera
$ws.Range("A98").Value = 1
$ws.Range("B98").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C98").Value = "Arica y Parinacota"
$ws.Range("D98").Value = 45009
$ws.Range("E98").Value = 15
$ws.Range("F98").Value = 100112021
$ws.Range("G98").Value = "Ají"
$ws.Range("H98").Value = "Inferno"
$ws.Range("I98").Value = "Primera"
$ws.Range("J98").Value = 190
$ws.Range("K98").Value = 24000
$ws.Range("L98").Value = 25000
$ws.Range("M98").Value = 24526
$ws.Range("N98").Value = "$/caja 15 kilos"
$ws.Range("O98").Value = "Región de Arica y Parinacota"
$ws.Range("P98").Value = 1635
$ws.Range("Q98").Value = 15
$ws.Range("R98").Value = "Hortaliza"

# Row 99: Inferno / Segunda
$ws.Range("A99").Value = 1
$ws.Range("B99").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C99").Value = "Arica y Parinacota"
$ws.Range("D99").Value = 45009
$ws.Range("E99").Value = 15
$ws.Range("F99").Value = 100112021
$ws.Range("G99").Value = "Ají"
$ws.Range("H99").Value = "Inferno"
$ws.Range("I99").Value = "Segunda"
$ws.Range("J99").Value = 135
$ws.Range("K99").Value = 21000
$ws.Range("L99").Value = 22000
$ws.Range("M99").Value = 21593
$ws.Range("N99").Value = "$/caja 15 kilos"
$ws.Range("O99").Value = "Región de Arica y Parinacota"
$ws.Range("P99").Value = 1440
$ws.Range("Q99").Value = 15
$ws.Range("R99").Value = "Hortaliza"
